$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the percentage figures up from row 3 to row 2 (E2:G2),
# matching the E3:G3 number format (0%) so the same style is reused.
$ws.Range("E2:G2").NumberFormat = $ws.Range("E3:G3").NumberFormat
$ws.Range("E2").Value = 0.05
$ws.Range("F2").Value = 0.05
$ws.Range("G2").Value = 0.9

# Row 3's percentage cells keep their formatting but lose their values.
$ws.Range("E3:G3").ClearContents()

# Update the active selection to reflect where the user ended up.
$ws.Range("G9").Select()
